$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"1.046867666666667"
$ws.Range("H2").Value = [double]"3.140603"
$ws.Range("I2").Value = [double]"0.000687505225377314"
$ws.Range("J2").Value = [double]"0.000687505225377314"
$ws.Range("M2").Value = [double]"0.141694"
$ws.Range("N2").Value = [double]"0.425082"
$ws.Range("O2").Value = [double]"0.01763793963212447"
$ws.Range("P2").Value = [double]"0.01763793963212447"
$ws.Range("Q2").Value = [double]"0.1483348671606667"
$ws.Range("R2").Value = [double]"1.335013804446"
$ws.Range("S2").Value = [double]"1.212617566197519E-05"
$ws.Range("T2").Value = [double]"1.212617566197519E-05"
$ws.Range("G3").Value = [double]"1.046867666666667"
$ws.Range("H3").Value = [double]"3.140603"
$ws.Range("I3").Value = [double]"0.000687505225377314"
$ws.Range("J3").Value = [double]"0.000687505225377314"
$ws.Range("O3").Value = [double]"0.2714637835982539"
$ws.Range("P3").Value = [double]"0.2714637835982538"
$ws.Range("Q3").Value = [double]"2.283007262687222"
$ws.Range("R3").Value = [double]"20.547065364185"
$ws.Range("S3").Value = [double]"0.0001866327697244959"
$ws.Range("T3").Value = [double]"0.0001866327697244959"
$ws.Range("G4").Value = [double]"1.046867666666667"
$ws.Range("H4").Value = [double]"3.140603"
$ws.Range("I4").Value = [double]"0.000687505225377314"
$ws.Range("J4").Value = [double]"0.000687505225377314"
$ws.Range("M4").Value = [double]"5.710985666666667"
$ws.Range("N4").Value = [double]"17.132957"
$ws.Range("O4").Value = [double]"0.7108982767696218"
$ws.Range("P4").Value = [double]"0.7108982767696217"
$ws.Range("Q4").Value = [double]"5.978646239230112"
$ws.Range("R4").Value = [double]"53.807816153071"
$ws.Range("S4").Value = [double]"0.000488746279990843"
$ws.Range("T4").Value = [double]"0.0004887462799908429"
$ws.Range("H5").Value = [double]"4442.55542"
$ws.Range("I5").Value = [double]"0.9725138978974124"
$ws.Range("J5").Value = [double]"0.9725138978974125"
$ws.Range("M5").Value = [double]"0.141694"
$ws.Range("N5").Value = [double]"0.425082"
$ws.Range("O5").Value = [double]"0.01763793963212447"
$ws.Range("P5").Value = [double]"0.01763793963212447"
$ws.Range("Q5").Value = [double]"209.8278158938266"
$ws.Range("R5").Value = [double]"1888.45034304444"
$ws.Range("S5").Value = [double]"0.01715314142251662"
$ws.Range("T5").Value = [double]"0.01715314142251662"
$ws.Range("H6").Value = [double]"4442.55542"
$ws.Range("I6").Value = [double]"0.9725138978974124"
$ws.Range("J6").Value = [double]"0.9725138978974125"
$ws.Range("O6").Value = [double]"0.2714637835982539"
$ws.Range("P6").Value = [double]"0.2714637835982538"
$ws.Range("S6").Value = [double]"0.2640023023251176"
$ws.Range("T6").Value = [double]"0.2640023023251175"
$ws.Range("H7").Value = [double]"4442.55542"
$ws.Range("I7").Value = [double]"0.9725138978974124"
$ws.Range("J7").Value = [double]"0.9725138978974125"
$ws.Range("M7").Value = [double]"5.710985666666667"
$ws.Range("N7").Value = [double]"17.132957"
$ws.Range("O7").Value = [double]"0.7108982767696218"
$ws.Range("P7").Value = [double]"0.7108982767696217"
$ws.Range("Q7").Value = [double]"8457.123442330771"
$ws.Range("R7").Value = [double]"76114.11098097694"
$ws.Range("S7").Value = [double]"0.6913584541497784"
$ws.Range("T7").Value = [double]"0.6913584541497784"
$ws.Range("G8").Value = [double]"40.80635833333333"
$ws.Range("H8").Value = [double]"122.419075"
$ws.Range("I8").Value = [double]"0.02679859687721029"
$ws.Range("J8").Value = [double]"0.0267985968772103"
$ws.Range("M8").Value = [double]"0.141694"
$ws.Range("N8").Value = [double]"0.425082"
$ws.Range("O8").Value = [double]"0.01763793963212447"
$ws.Range("P8").Value = [double]"0.01763793963212447"
$ws.Range("Q8").Value = [double]"5.782016137683332"
$ws.Range("R8").Value = [double]"52.03814523914999"
$ws.Range("S8").Value = [double]"0.0004726720339458745"
$ws.Range("T8").Value = [double]"0.0004726720339458745"
$ws.Range("G9").Value = [double]"40.80635833333333"
$ws.Range("H9").Value = [double]"122.419075"
$ws.Range("I9").Value = [double]"0.02679859687721029"
$ws.Range("J9").Value = [double]"0.0267985968772103"
$ws.Range("O9").Value = [double]"0.2714637835982539"
$ws.Range("P9").Value = [double]"0.2714637835982538"
$ws.Range("Q9").Value = [double]"88.99043824273609"
$ws.Range("R9").Value = [double]"800.9139441846249"
$ws.Range("S9").Value = [double]"0.007274848503411859"
$ws.Range("T9").Value = [double]"0.007274848503411858"
$ws.Range("G10").Value = [double]"40.80635833333333"
$ws.Range("H10").Value = [double]"122.419075"
$ws.Range("I10").Value = [double]"0.02679859687721029"
$ws.Range("J10").Value = [double]"0.0267985968772103"
$ws.Range("M10").Value = [double]"5.710985666666667"
$ws.Range("N10").Value = [double]"17.132957"
$ws.Range("O10").Value = [double]"0.7108982767696218"
$ws.Range("P10").Value = [double]"0.7108982767696217"
$ws.Range("Q10").Value = [double]"233.0445275505305"
$ws.Range("R10").Value = [double]"2097.400747954775"
$ws.Range("S10").Value = [double]"0.01905107633985257"
$ws.Range("T10").Value = [double]"0.01905107633985257"
